$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.891.09'
$ws.Range('E2').Value = '  -1.35%  '
$ws.Range('D3').Value = '3.158.99'
$ws.Range('E3').Value = '  -4.77%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '590.72'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.17'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.77%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '3.158.27'
$ws.Range('E8').Value = '  -4.79%  '
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('E10').Value = '  -5.51%  '
$ws.Range('E11').Value = '  -5.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.453'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.31%  '
$ws.Range('E13').Value = '  -4.20%  '
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('D15').Value = '3.679.47'
$ws.Range('E15').Value = '  -4.72%  '
$ws.Range('E16').Value = '  -1.85%  '
$ws.Range('D17').Value = '3.164.23'
$ws.Range('E17').Value = '  -4.47%  '
$ws.Range('D18').Value = '62.841.28'
$ws.Range('E18').Value = '  -1.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '460.48'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.09%  '
$ws.Range('E21').Value = '  -1.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.698'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.60'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.35'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.47'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.48%  '
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('E28').Value = '  -3.84%  '
$ws.Range('E29').Value = '  -6.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.72'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.84%  '
$ws.Range('E31').Value = '  -6.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.08'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.35%  '
$ws.Range('E33').Value = '  -2.68%  '
$ws.Range('E34').Value = '  -6.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.04'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.46%  '
$ws.Range('E36').Value = '  -4.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.10'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.21%  '
$ws.Range('D38').Value = '0.0₃0702'
$ws.Range('E38').Value = '  -5.37%  '
$ws.Range('E39').Value = '  -2.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '400.98'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.11'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.57%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.112'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.01%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.61'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.53%  '
$ws.Range('D44').Value = '2.799.26'
$ws.Range('E44').Value = '  -9.65%  '
$ws.Range('E45').Value = '  -5.17%  '
$ws.Range('E47').Value = '  -5.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '25.32'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.44'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('E50').Value = '  -2.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '34.17'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.96%  '
